# Generate Report for Handoff
# Adds two new handoff entries (a .png file and a .md file) to the
# localization-status report: the Overview sheet gets two more rows,
# and each language sheet (zh-cn / de-de) gets matching detail rows.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# New "source" file identities being reported on in this handoff batch.
# ---------------------------------------------------------------------
$png1Name = "76b2e1a1-7b18-44ac-9c75-0fc145233969.png"
$png2Name = "87bb59a7-c75f-4ef2-982d-ddb3d676f8e8.png"
$mdName   = "de0cd63a-60ab-41dc-8236-f71a1143109b.md"

$png1TargetZh = "110d03b8f05d1bb5b3d3c855bfbb37ae9b54d352.png"
$png2TargetZh = "80d35c54a121d1a14ed9a8c8f965e770aa61833a.png"
$mdTargetZh   = "de0cd63a-60ab-41dc-8236-f71a1143109b.db86cac4270845698a910eff5e9cf56c290c62e6.zh-cn.xlf"

$png1TargetDe = $png1TargetZh
$png2TargetDe = $png2TargetZh
$mdTargetDe   = "de0cd63a-60ab-41dc-8236-f71a1143109b.db86cac4270845698a910eff5e9cf56c290c62e6.de-de.xlf"

$png1Url = "https://github.com/OpenLocalizationTest/oltest/blob/2dd0f1cea29b39ed8fe6f49597b0628dcd430fea/e2e/$png1Name"
$png2Url = "https://github.com/OpenLocalizationTest/oltest/blob/2dd0f1cea29b39ed8fe6f49597b0628dcd430fea/e2e/$png2Name"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/2dd0f1cea29b39ed8fe6f49597b0628dcd430fea/e2e/$mdName"

$png1TargetZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/652d433928a04fa7f884047e990711f7fc3f67ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1TargetZh"
$png2TargetZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/652d433928a04fa7f884047e990711f7fc3f67ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2TargetZh"
$mdTargetZhUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/652d433928a04fa7f884047e990711f7fc3f67ea/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$mdTargetZh"

$png1TargetDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69f91f418d46dc4a17ff4dce032468f4af64fefb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1TargetDe"
$png2TargetDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69f91f418d46dc4a17ff4dce032468f4af64fefb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2TargetDe"
$mdTargetDeUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/69f91f418d46dc4a17ff4dce032468f4af64fefb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$mdTargetDe"

$handoffDateTime = "2016-49-18 14:49:11"
$zhHandoffDt     = "2016-03-18 14:49:09"
$deHandoffDt     = "2016-03-18 14:49:11"
$epoch           = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet: update row 2, add rows 3 and 4
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Hyperlinks.Delete() | Out-Null

$wsOverview.Cells.Item(2, 1).Value = $png1Name
$wsOverview.Cells.Item(2, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(2, 4).Value = $handoffDateTime

$wsOverview.Cells.Item(3, 1).Value = $png2Name
$wsOverview.Cells.Item(3, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 4).Value = $handoffDateTime

$wsOverview.Cells.Item(4, 1).Value = $mdName
$wsOverview.Cells.Item(4, 2).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 3).Value = "Ready for handoff"
$wsOverview.Cells.Item(4, 4).Value = $handoffDateTime

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $png1Url, "", "", $png1Name) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $png2Url, "", "", $png2Name) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $mdUrl, "", "", $mdName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet: add rows 3 and 4 (row 2 is refreshed with the same
# source file but a new target/handoff file)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete() | Out-Null
$wsZh.Range("B2").Hyperlinks.Delete() | Out-Null
$wsZh.Range("D2").Hyperlinks.Delete() | Out-Null

$wsZh.Cells.Item(2, 1).Value = $png1Name
$wsZh.Cells.Item(2, 2).Value = ".png"
$wsZh.Cells.Item(2, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(2, 4).Value = $png1TargetZh
$wsZh.Cells.Item(2, 5).Value = $zhHandoffDt
$wsZh.Cells.Item(2, 5).NumberFormat = $dateFmt
$wsZh.Cells.Item(2, 8).Value = $epoch
$wsZh.Cells.Item(2, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item(2, 9).Value = "IsDependency"
$wsZh.Cells.Item(2, 10).Value = "e2e\$mdName"

$wsZh.Cells.Item(3, 1).Value = $png2Name
$wsZh.Cells.Item(3, 2).Value = ".png"
$wsZh.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(3, 4).Value = $png2TargetZh
$wsZh.Cells.Item(3, 5).Value = $zhHandoffDt
$wsZh.Cells.Item(3, 5).NumberFormat = $dateFmt
$wsZh.Cells.Item(3, 8).Value = $epoch
$wsZh.Cells.Item(3, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item(3, 9).Value = "IsDependency"
$wsZh.Cells.Item(3, 10).Value = "e2e\$mdName"

$wsZh.Cells.Item(4, 1).Value = $mdName
$wsZh.Cells.Item(4, 2).Value = ".md"
$wsZh.Cells.Item(4, 3).Value = "Ready for handoff"
$wsZh.Cells.Item(4, 4).Value = $mdTargetZh
$wsZh.Cells.Item(4, 5).Value = $zhHandoffDt
$wsZh.Cells.Item(4, 5).NumberFormat = $dateFmt
$wsZh.Cells.Item(4, 8).Value = $epoch
$wsZh.Cells.Item(4, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item(4, 9).Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $png1Url, "", "", $png1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $png1Url, "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $png1TargetZhUrl, "", "", $png1TargetZh) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $png2Url, "", "", $png2Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $png2Url, "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $png2TargetZhUrl, "", "", $png2TargetZh) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $mdUrl, "", "", $mdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), $mdUrl, "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), $mdTargetZhUrl, "", "", $mdTargetZh) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet: same shape as zh-cn, but with the de-de handoff
# timestamp and de-de target/handback files
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete() | Out-Null
$wsDe.Range("B2").Hyperlinks.Delete() | Out-Null
$wsDe.Range("D2").Hyperlinks.Delete() | Out-Null

$wsDe.Cells.Item(2, 1).Value = $png1Name
$wsDe.Cells.Item(2, 2).Value = ".png"
$wsDe.Cells.Item(2, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(2, 4).Value = $png1TargetDe
$wsDe.Cells.Item(2, 5).Value = $deHandoffDt
$wsDe.Cells.Item(2, 5).NumberFormat = $dateFmt
$wsDe.Cells.Item(2, 8).Value = $epoch
$wsDe.Cells.Item(2, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item(2, 9).Value = "IsDependency"
$wsDe.Cells.Item(2, 10).Value = "e2e\$mdName"

$wsDe.Cells.Item(3, 1).Value = $png2Name
$wsDe.Cells.Item(3, 2).Value = ".png"
$wsDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(3, 4).Value = $png2TargetDe
$wsDe.Cells.Item(3, 5).Value = $deHandoffDt
$wsDe.Cells.Item(3, 5).NumberFormat = $dateFmt
$wsDe.Cells.Item(3, 8).Value = $epoch
$wsDe.Cells.Item(3, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item(3, 9).Value = "IsDependency"
$wsDe.Cells.Item(3, 10).Value = "e2e\$mdName"

$wsDe.Cells.Item(4, 1).Value = $mdName
$wsDe.Cells.Item(4, 2).Value = ".md"
$wsDe.Cells.Item(4, 3).Value = "Ready for handoff"
$wsDe.Cells.Item(4, 4).Value = $mdTargetDe
$wsDe.Cells.Item(4, 5).Value = $deHandoffDt
$wsDe.Cells.Item(4, 5).NumberFormat = $dateFmt
$wsDe.Cells.Item(4, 8).Value = $epoch
$wsDe.Cells.Item(4, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item(4, 9).Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $png1Url, "", "", $png1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $png1Url, "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $png1TargetDeUrl, "", "", $png1TargetDe) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $png2Url, "", "", $png2Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $png2Url, "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $png2TargetDeUrl, "", "", $png2TargetDe) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $mdUrl, "", "", $mdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), $mdUrl, "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), $mdTargetDeUrl, "", "", $mdTargetDe) | Out-Null
